# Aula 29 - Incluindo o Thymeleaf-layout - dependencia e alteracao do application.properties
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$obs1 = "5:56 - abordado chaves e configurações que podem ser inseridas no application.properties para mais configurações do thymeleaf. O próprio spring boot ja pré-configura algumas coisas do thymeleaf, mas caso não seja usado o spring boot, podem ser feitas estas alterações direto no projeto, como por exemplo o prefixo de paginas citado em aulas anteriores, que é o diretorio onde fica localizado os arquivos .html do front-end. Mais detalhes na documentação descrita no link: https://docs.spring.io/spring-boot/docs/current/reference/htmlsingle/"
$sessao = "7. Thymeleaf para as Views"
$nomeAula = "29. Incluindo o Thymeleaf-Layout"
$obs2 = "5:03`ninclusão de dependencia para uso de templates no thymeleaf`n<dependency>`n<groupId>nz.net.ultraq.thymeleaf</groupId>`n<artifactId>thymeleaf-layout-dialect</artifactId>`n</dependency>"

# Row 17
$ws.Cells.Item(17, 2).Value = 29
$ws.Cells.Item(17, 5).Value = $obs1
$ws.Cells.Item(17, 3).Value = $sessao
$ws.Cells.Item(17, 4).Value = $nomeAula
$ws.Cells.Item(17, 5).WrapText = $true
$ws.Rows.Item(17).RowHeight = 120

# Row 18
$ws.Cells.Item(18, 2).Value = 29
$ws.Cells.Item(18, 3).Value = $sessao
$ws.Cells.Item(18, 4).Value = $nomeAula
$ws.Cells.Item(18, 5).Value = $obs2
$ws.Cells.Item(18, 5).WrapText = $true
$ws.Rows.Item(18).RowHeight = 105

# Update view/selection to match the post-edit state
$ws.Application.ActiveWindow.ScrollRow = 16
$ws.Range("E21").Select()
